# Fill in the real figure captions (replacing the placeholder "Рис. N: РисN"
# text) for the 5 screenshots in the lab-3 report, both in the picture's
# alt text (wp:docPr/@descr, exposed as InlineShape.AlternativeText) and in
# the visible "Image Caption" paragraph that follows each figure.

$d = $word.ActiveDocument

$map = @{
    "Рис. 1: Рис1" = "Рис. 1: Переход в каталог курса"
    "Рис. 2: Рис2" = "Рис. 2: Компиляция шаблона"
    "Рис. 3: Рис3" = "Рис. 3: Удаление файла"
    "Рис. 4: Рис4" = "Рис. 4: Редактируем файл report.md"
    "Рис. 5: Рис5" = "Рис. 5: Загрузка файлов на github"
}

# --- 1) Update each picture's alternative text (wp:docPr descr="...") ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $old = $shape.AlternativeText
    if ($map.ContainsKey($old)) {
        $shape.AlternativeText = $map[$old]
    }
}

# --- 2) Update the visible caption paragraphs (style "Image Caption") ---
foreach ($para in $d.Paragraphs) {
    $style = $para.Style
    if ($style -ne $null -and $style.NameLocal -eq "Image Caption") {
        $old = $para.Range.Text
        foreach ($key in $map.Keys) {
            if ($old.Contains($key)) {
                $para.Range.Text = $map[$key]
            }
        }
    }
}
